# Helper: find `searchText` starting the search at character offset `startPos`
# (searches from startPos to the end of the document body). Returns the found
# Range (collapsed to the match) or $null if not found.
function Find-FromPos {
    param($doc, $searchText, $startPos)
    $endPos = $doc.Content.End
    $r = $doc.Range($startPos, $endPos)
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { return $null }
    return $r
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Purpose of the Functional Safety Concept" heading: the original run is
#    split as "Purpose of th" + (_GoBack bookmark) + "e Functional Safety
#    Concept". Re-typing the full heading as one Find/Replace merges the runs
#    back into a single run and drops the (now-empty) _GoBack bookmark, which
#    also renumbers every later bookmark id down by one - matching the diff.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Purpose of the Functional Safety Concept", $false, $false, $false, $false,
    $false, $true, 1, $false, "Purpose of the Functional Safety Concept", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge the two runs " the general functionality of the item" + "."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The Functional Safety Concept document the general functionality of the item.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The Functional Safety Concept document the general functionality of the item.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Description of architecture elements" heading: "D" + "escription..." ->
#    single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Description of architecture elements", $false, $false, $false, $false,
    $false, $true, 1, $false, "Description of architecture elements", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Add <w:noProof/> to the rPr of the runs that hold the two non-cover
#    inline pictures (Preliminary Architecture diagram + Refinement of the
#    System Architecture diagram). InlineShapes #1 is the cover picture
#    (already has noProof); #2 and #3 are the ones touched by the diff.
# ---------------------------------------------------------------------------
$d.InlineShapes.Item(2).Range.NoProofing = $true
$d.InlineShapes.Item(3).Range.NoProofing = $true

# ---------------------------------------------------------------------------
# 5) "Allocation of Functional Safety Requirements to Architecture Elements"
#    heading: "Alloc" + "ation..." -> single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Allocation of Functional Safety Requirements to Architecture Elements",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Allocation of Functional Safety Requirements to Architecture Elements", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Functional Safety Requirements table (first/"LDW Requirements" table):
#    row 1 (amplitude) - wording + typo fix, no bookmark here.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The lane keeping item shall ensure that the lane departure oscillating torque amplitude is below ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The EPS ECU shall ensure that the lane departure warning torque amplitude is below ", 2) | Out-Null

$d.Content.Find.Execute(
    "Max_Torque_Amplitude", $false, $false, $false, $false, $false, $true, 1, $false,
    "Max_Torque_Ampliture", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Same table, row 2 (frequency) - wording fix only.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The lane keeping item shall ensure that the lane departure oscillating torque frequency is below ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The EPS ECU shall ensure that the lane departure warning torque frequency is below ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8) Verification table: "Max_Torque_" + "Frequency" -> single run
#    "Max_Torque_Frequency" (the sentence text around it is unchanged there).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Verify that the system does turn off in time if Max_Torque_Frequency is exceeded.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Verify that the system does turn off in time if Max_Torque_Frequency is exceeded.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 9) Allocation table header cell: "Functi" + "onal Safety Requirement" ->
#    single run. (Two other "Functional Safety Requirement" header cells are
#    already single runs and Find only matches the still-split one.)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Functional Safety Requirement", $false, $false, $false, $false, $false,
    $true, 1, $false, "Functional Safety Requirement", 2) | Out-Null

# ---------------------------------------------------------------------------
# 10) Allocation table, row 1 (amplitude): wording + typo fix, and wrap the
#     paragraph content with a new _GoBack bookmark (start right at the top
#     of the paragraph, end right after the "Max_Torque_Ampliture" run).
# ---------------------------------------------------------------------------
$sentenceOld = "The lane keeping item shall ensure that the lane departure oscillating torque amplitude is below "
$sentenceStart = Find-FromPos $d $sentenceOld 0
$paraStart = $sentenceStart.Start

$d.Content.Find.Execute(
    $sentenceOld, $false, $false, $false, $false, $false, $true, 1, $false,
    "The EPS ECU shall ensure that the lane departure warning torque amplitude is below ", 2) | Out-Null

$amplitudeRun = Find-FromPos $d "Max_Torque_Amplitude" 0
$d.Content.Find.Execute(
    "Max_Torque_Amplitude", $false, $false, $false, $false, $false, $true, 1, $false,
    "Max_Torque_Ampliture", 2) | Out-Null

$bmEnd = $paraStart + ("The EPS ECU shall ensure that the lane departure warning torque amplitude is below Max_Torque_Ampliture").Length
$bmRange = $d.Range($paraStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 11) Allocation table, row 2 (frequency): wording fix only.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The lane keeping item shall ensure that the lane departure oscillating torque frequency is below ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The EPS ECU shall ensure that the lane departure warning torque frequency is below ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 12) Warning/Degradation table: "Malfunction_0" + "2" -> single run, and
#     "Malfunction_0" + "3" -> single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Malfunction_02", $false, $false, $false, $false, $false, $true, 1, $false,
    "Malfunction_02", 2) | Out-Null

$d.Content.Find.Execute(
    "Malfunction_03", $false, $false, $false, $false, $false, $true, 1, $false,
    "Malfunction_03", 2) | Out-Null

Write-Output "done"
